$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155-179 down to 156-180.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new data point
# (same static/category fields as its neighbours, new date + volume).
$ws.Cells.Item(155, 1).Value = 10
$ws.Cells.Item(155, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(155, 3).Value = "La Araucanía"
$ws.Cells.Item(155, 4).Value = "11/22/2021"
$ws.Cells.Item(155, 5).Value = 9
$ws.Cells.Item(155, 6).Value = 100112039
$ws.Cells.Item(155, 7).Value = "Ciboulette"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 30
$ws.Cells.Item(155, 11).Value = 5000
$ws.Cells.Item(155, 12).Value = 5000
$ws.Cells.Item(155, 13).Value = 5000
$ws.Cells.Item(155, 14).Value = "$/docena de atados"
$ws.Cells.Item(155, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(155, 16).Value = 1667
$ws.Cells.Item(155, 17).Value = 3
$ws.Cells.Item(155, 18).Value = "Hortaliza"
